# 7_gunluk_vardiya_plani.xlsx - data refresh for the shift-plan sheet.
# Updates the hourly headcount figures in rows 4-10 (Pazartesi..Pazar) and
# normalizes the header row (B1:Y1) back onto the same cell style used by
# A1, collapsing the redundant duplicate style entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row style cleanup: B1:Y1 should share A1's style (border +
# bold + centered/top aligned), instead of a duplicate style index. ---
$hdr = $ws.Range("B1:Y1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# --- Pazartesi (row 4) ---
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 4
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 6
$ws.Range("R4").Value = 8
$ws.Range("S4").Value = 8
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 8
$ws.Range("V4").Value = 7
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 1

# --- Sali (row 5) ---
$ws.Range("N5").Value = 5
$ws.Range("P5").Value = 6
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 6
$ws.Range("U5").Value = 7
$ws.Range("X5").Value = 1

# --- Carsamba (row 6) ---
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 7
$ws.Range("O6").Value = 8
$ws.Range("P6").Value = 8
$ws.Range("S6").Value = 9
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = 9
$ws.Range("W6").Value = 6
$ws.Range("X6").Value = 2

# --- Persembe (row 7) ---
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 6
$ws.Range("O7").Value = 7
$ws.Range("R7").Value = 9
$ws.Range("S7").Value = 8
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 7
$ws.Range("X7").Value = 2

# --- Cuma (row 8) ---
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 3
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = 7
$ws.Range("P8").Value = 7
$ws.Range("Q8").Value = 7
$ws.Range("R8").Value = 7
$ws.Range("S8").Value = 7
$ws.Range("U8").Value = 8
$ws.Range("X8").Value = 1

# --- Cumartesi (row 9) ---
$ws.Range("L9").Value = 2
$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 6
$ws.Range("P9").Value = 6
$ws.Range("Q9").Value = 7
$ws.Range("R9").Value = 7
$ws.Range("S9").Value = 7
$ws.Range("T9").Value = 7
$ws.Range("V9").Value = 7
$ws.Range("X9").Value = 1

# --- Pazar (row 10) ---
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 3
$ws.Range("O10").Value = 6
$ws.Range("P10").Value = 6
$ws.Range("Q10").Value = 6
$ws.Range("R10").Value = 7
$ws.Range("S10").Value = 7
$ws.Range("T10").Value = 7
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 1
